# Commit: "Added Test Data For Hungary/Russia/Finland Market"
#
# Adds three new worksheets (Russia, Finland, Hungary) after "Denmark",
# each cloned from the "Denmark" sheet so the layout/styles/merged cells
# match the existing per-country test-data sheets, then fills in the
# country-specific "<Country> Market" label (B2) and Jira ticket (B4).
#
# NOTE: this COM-interop runtime only binds *positional* parameters on
# user-defined functions, so per-sheet config is driven from an array of
# hashtables + foreach instead of a parameterised helper function.

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("Denmark")

$newMarkets = @(
    @{ Name = "Russia";  Ticket = "NGC-2929/T2898"; Selection = "A1:D11" },
    @{ Name = "Finland"; Ticket = "NGC-3130/T2941"; Selection = "A1:D11" },
    @{ Name = "Hungary"; Ticket = "NGC-3104/T2990"; Selection = "K11" }
)

foreach ($market in $newMarkets) {
    # Clone the Denmark template sheet and drop the copy at the end of the
    # tab strip (After:= last worksheet).
    $template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $market.Name

    # B4 = ticket reference, B2 = "<Country> Market" label. Ticket first so
    # the shared-string insertion order matches the source edit.
    $newSheet.Range("B4").Value = $market.Ticket
    $newSheet.Range("B2").Value = "$($market.Name) Market"

    # Make each newly-created sheet active as it is added (mirrors a user
    # clicking through the new tabs while filling them in) and restore its
    # on-screen selection.
    $newSheet.Activate()
    $newSheet.Range($market.Selection).Select()
}
